# fix: unique command names in XLSX - prefix protocol name to each step
#
# For every "protocol" worksheet (price1, price2, discount1, ... boosters),
# prefix each Step/command name in column A (rows below the header row)
# with the sheet's own name, e.g. "Step4 Seed" -> "free1 Step4 Seed".

$wb = $excel.ActiveWorkbook

$targetSheets = @(
    "price1","price2","discount1","discount2","free1","free2",
    "nomoney1","nomoney2","noppv1","noppv2","card1","card2",
    "nosex1","nosex2","offtopic1","offtopic2","real1","real2",
    "voice1","voice2","customyes1","customyes2","customno1","customno2",
    "done1","done2","cumcontrol","dickpic","boosters"
)

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    $lastRow = $used.Row + $used.Rows.Count - 1

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $val = $cell.Value()
        if ($val -ne $null -and $val -ne "") {
            $prefix = $sheetName + " "
            if ($val.ToString().StartsWith($prefix) -eq $false) {
                $cell.Value = $prefix + $val
            }
        }
    }
}
